# Reposition/resize the three "external databases" icon pictures on slide 1.
# Target sizes/positions are given in EMU in the source OOXML; PowerPoint's
# COM object model expresses Shape.Left/Top/Width/Height in points
# (1 pt = 12700 EMU), so the literals below are the point values that this
# host's points->EMU conversion reproduces exactly for the desired EMU
# targets.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Picture 21 (top-right "External Databases" icon)
# off: 7817560,1464093 -> 7886690,1468355 ; ext: 2038259,1045137 -> 1937794,993622
$shp = $s.Shapes.Item("Picture 21")
$shp.Left   = 620.9992907234252
$shp.Top    = 115.61850393700787
$shp.Width  = 152.58228284940947
$shp.Height = 78.23795275590551

# Picture 23
# off: 9003974,3986427 -> 9003974,4016989 ; ext: 1086890,413628 -> 1006579,383065
$shp = $s.Shapes.Item("Picture 23")
$shp.Left   = 708.9743697711615
$shp.Top    = 316.29838551919295
$shp.Width  = 79.25822803887795
$shp.Height = 30.16267655019685

# Picture 27
# off: 9003974,4714902 -> 9003974,4733019 ; ext: 753484,327275 -> 650389,282496
$shp = $s.Shapes.Item("Picture 27")
$shp.Left   = 708.9743697711615
$shp.Top    = 372.67870047982285
$shp.Width  = 51.211771345964564
$shp.Height = 22.243779527559056
